# Conserto do erro com o rótulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
#
# For each of the first five sheets (the ones with a B1:E1 year-header row),
# cell E1 incorrectly held a stray numeric value (705.7799333869282) instead
# of the "2050" (or "2041-2050") text label that belongs alongside the other
# year headers. We fix that label, then strip the "Total" row that closes
# out every table (rows 13 on sheets 1-4, row 4 on sheet 6 "Custo Total").
# Sheet 5 ("Emissoes Totais") only needs the label fix - it never had a
# Total row.

# Helper: write a text value into a cell via a scratch cell formatted as
# Text, then Copy/PasteSpecial(values) onto the target. Using .Value
# directly on a numeric-looking string ("2050") lets Excel auto-coerce it
# to a number, which is exactly the historical bug we're fixing - the
# round-trip through a Text-formatted scratch cell keeps the destination's
# original style (border/bold/center) while forcing a genuine text value.
function Set-TextLabel {
    param($Sheet, $Target, $Text)

    $scratch = $Sheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $Sheet.Range($Target).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Potencia Acumulada - SIN (MW) ---
$ws1 = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-TextLabel $ws1 "E1" "2050"
$ws1.Rows("13:13").Delete()

# --- Sheet 2: Geracao Periodo Medio (MWMed) ---
$ws2 = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-TextLabel $ws2 "E1" "2050"
$ws2.Rows("13:13").Delete()

# --- Sheet 3: Atendimento a Ponta(MW) ---
$ws3 = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-TextLabel $ws3 "E1" "2050"
$ws3.Rows("13:13").Delete()

# --- Sheet 4: Potencia Incremental - SIN(MW) (decade ranges, not years) ---
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextLabel $ws4 "E1" "2041-2050"
$ws4.Rows("13:13").Delete()

# --- Sheet 5: Emissoes Totais (MtCO2eq) - label fix only, no Total row ---
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabel $ws5 "E1" "2050"

# --- Sheet 6: Custo Total (bilhões de R$) ---
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows("4:4").Delete()
